$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2; this shifts existing rows 2-36 down to 3-37
# (values + per-cell styles move with them).
$ws.Rows("2:2").Insert()

# The freshly inserted row 2 picks up header-like formatting; copy the correct
# per-column formats from row 3 (which now holds what used to be row 2's data).
$ws.Range("A3:F3").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new row 2 with the latest circular entry.
$ws.Range("A2").Value = 36
$ws.Range("B2").Value = "2. P0610 (99.85% min) /P1020/ EC Grade Ingot & Sow 99.7% (min) / Cast Bar"
$ws.Range("C2").Value = "P1020"
$ws.Range("D2").Value = 265.5
$ws.Range("E2").Value = "25.09.2025"
$ws.Range("F2").Value = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-25-september-2025.pdf"

# Hyperlinks are anchored to fixed ranges and do not shift with the row
# insert, so rebuild the whole collection against the new row numbers.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-25-september-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-23-september-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-20-september-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-18-september-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-17-september-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-13-september-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-12-september-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-02-september-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-27-august-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F11"), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-26-august-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F12"), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-23-august-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F13"), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-20-august-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F14"), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-19-august-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F15"), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-14-august-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F16"), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-13-august-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F17"), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-12-august-2025.pdf")

# Adding a hyperlink auto-applies the blue/underlined "Hyperlink" style; the
# source sheet keeps the plain centered style on those cells, so restore it
# by pasting the format from an unaffected cell in the same column.
$ws.Range("F18").Copy()
$ws.Range("F2:F17").PasteSpecial(-4122)
$excel.CutCopyMode = 0
